# Update gh-pages output data (苏州-漫展信息.xlsx)
# Applies refreshed "想去人数" (col F) and "最低票价" (col G) values
# to the 展览 (sheet1) and 全部类型 (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 84
$ws1.Range("F3").Value  = 21807
$ws1.Range("F5").Value  = 8143
$ws1.Range("G6").Value  = 58
$ws1.Range("G8").Value  = 55
$ws1.Range("G9").Value  = 65
$ws1.Range("F11").Value = 549
$ws1.Range("F13").Value = 250
$ws1.Range("F14").Value = 820
$ws1.Range("F15").Value = 1375
$ws1.Range("F18").Value = 723
$ws1.Range("F20").Value = 108
$ws1.Range("F21").Value = 101
$ws1.Range("F23").Value = 1242
$ws1.Range("F24").Value = 85
$ws1.Range("F25").Value = 57
$ws1.Range("F26").Value = 248
$ws1.Range("F27").Value = 5213
$ws1.Range("F28").Value = 625
$ws1.Range("F31").Value = 5289
$ws1.Range("F34").Value = 77
$ws1.Range("F36").Value = 13567
$ws1.Range("F37").Value = 1390
$ws1.Range("F38").Value = 167
$ws1.Range("F39").Value = 67
$ws1.Range("F41").Value = 357
$ws1.Range("F42").Value = 493
$ws1.Range("F43").Value = 4110
$ws1.Range("F44").Value = 55
$ws1.Range("F45").Value = 338
$ws1.Range("F46").Value = 106

# --- Sheet "全部类型" ---------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 84
$ws4.Range("F3").Value  = 21807
$ws4.Range("F5").Value  = 8143
$ws4.Range("G6").Value  = 58
$ws4.Range("G8").Value  = 55
$ws4.Range("G9").Value  = 65
$ws4.Range("F11").Value = 549
$ws4.Range("F13").Value = 250
$ws4.Range("F14").Value = 822
$ws4.Range("F15").Value = 1375
$ws4.Range("F18").Value = 723
$ws4.Range("F20").Value = 108
$ws4.Range("F21").Value = 101
$ws4.Range("F23").Value = 1242
$ws4.Range("F24").Value = 85
$ws4.Range("F25").Value = 57
$ws4.Range("F26").Value = 248
$ws4.Range("F28").Value = 5213
$ws4.Range("F29").Value = 625
$ws4.Range("F34").Value = 5289
$ws4.Range("F37").Value = 77
$ws4.Range("F39").Value = 13567
$ws4.Range("F40").Value = 1390
$ws4.Range("F41").Value = 167
$ws4.Range("F42").Value = 67
$ws4.Range("F44").Value = 357
$ws4.Range("F45").Value = 493
$ws4.Range("F46").Value = 4110
$ws4.Range("F47").Value = 55
$ws4.Range("F48").Value = 338
$ws4.Range("F49").Value = 106
